# BIS-772: Add "Organism" as one more selected field (PropertyTypePermId)
# to the XLS export header row, right after "Comment", to improve the
# sorting check.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 4) currently ends at column L ("Comment").
# Append a new header cell "Organism" in column M, matching the
# formatting of the rest of the header row.
$ws.Range("M4").Value = "Organism"

# Copy the formatting of the last existing header cell onto the new one
# so the new header cell looks consistent with its neighbours.
$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null

# Move/update the active selection to the newly added header cell,
# mirroring the selection change recorded for this edit.
$ws.Range("M4").Select() | Out-Null
